$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header for the "19-jul" column in AG1
$ws.Range("AG1").Value = "19-jul"

# Values for the new AG column (rows 2-18)
$values = @(0, 10, 24, 26, 0, 9, 10, 21, 23, 12, 0, 7, 0, 0, 3, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 33).Value = $values[$i]
}

# Update selection to match the saved file (AG19)
$ws.Range("AG19").Select()
